$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new error-code rows (write order matches the original
# authoring sequence so new shared-string indices line up with the diff)
$ws.Range("A6").Value = "NO_BLOOD_TYPE_DATE"
$ws.Range("A7").Value = "NO_BLOOD_TYPE_RESULT"
$ws.Range("B7").Value = "No documented blood type result. Verify in patient's physical record and Genesis."
$ws.Range("C6").Value = "Blood Type Verification"
$ws.Range("C7").Value = "Blood Type Verification"
$ws.Range("B6").Value = "No documented blood type date (G6PD Date). Verify in patient's physical record and Genesis."

# Re-fit column widths to the new (longer) content, matching Excel's AutoFit.
# (The host quantizes ColumnWidth to whole-pixel/MDW-7 steps, so these inputs
# are chosen to land on the grid point nearest the recorded target width.)
$ws.Columns.Item(2).ColumnWidth = 78.66666666666667
$ws.Columns.Item(3).ColumnWidth = 19.333333333333332

# Match the saved selection state
$ws.Range("A7").Select()
